# Rename 'Codelists' sheet to 'Cells' and update the active sheet/selection
# to match the saved workbook view state (Close #256).

$wb = $excel.ActiveWorkbook

# Rename the "Codelists" sheet to "Cells"
$wsCells = $wb.Worksheets.Item("Codelists")
$wsCells.Name = "Cells"

# Make the "Cells" sheet the active/selected tab, with G13 selected
$wsCells.Activate()
$wsCells.Range("G13").Select()
